$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text number format to column A (so it displays/stores values as text)
$ws.Columns.Item(1).NumberFormat = "@"

# A3 becomes a text phone number with country code, normalized
$ws.Range("A3").Value = "+233558059073"

# Page orientation explicitly set to portrait
$ws.PageSetup.Orientation = 1

# Update active cell selection to D4 (cosmetic, matches diff)
$ws.Range("D4").Select()
